# Sprint 1 Architecture edit
# Applies:
#  1) Collapse the split "April " / "19, " / "2017" runs into a single
#     "April 19, 2017" run on every date textbox across the deck.
#  2) Collapse the split "Denial of " / "service" / ". " runs (slide 11)
#     into a single "Denial of service. " run.
#  3) Rewrite the Effort Estimate "Productivity (PROD)" slide (slide 9)
#     content placeholder with the new sprint numbers / team breakdown.

$p = $ppt.ActivePresentation

function Set-ParaText {
    param($shape, [int]$paraIndex, [string]$text)
    $tr = $shape.TextFrame.TextRange
    $para = $tr.Paragraphs($paraIndex, 1)
    # Force a full defragment of the paragraph's runs: touch it with a
    # throw-away value first so that when we set the real text back,
    # PowerPoint collapses everything into a single run sharing the
    # first run's formatting (matching how the authors' edit merged
    # runs without leaving stray empty runs behind).
    $para.Text = "~"
    $para.Text = $text
}

# ---- 1) Date textbox merges -------------------------------------------------
$dateTargets = @(
    @{ Slide = 2;  Shape = 6 },
    @{ Slide = 3;  Shape = 4 },
    @{ Slide = 4;  Shape = 5 },
    @{ Slide = 5;  Shape = 5 },
    @{ Slide = 6;  Shape = 5 },
    @{ Slide = 7;  Shape = 5 },
    @{ Slide = 8;  Shape = 5 },
    @{ Slide = 9;  Shape = 5 },
    @{ Slide = 10; Shape = 4 },
    @{ Slide = 11; Shape = 4 },
    @{ Slide = 13; Shape = 3 }
)

foreach ($t in $dateTargets) {
    $s = $p.Slides.Item($t.Slide)
    $shp = $s.Shapes.Item($t.Shape)
    Set-ParaText -shape $shp -paraIndex 1 -text "April 19, 2017"
}

# ---- 2) Slide 11: "Denial of service." run merge ---------------------------
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(2)
Set-ParaText -shape $shp11 -paraIndex 2 -text "Denial of service. "

# ---- 3) Slide 9: Productivity (PROD) content rewrite ------------------------
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$tr9 = $shp9.TextFrame.TextRange

$lines = @(
    "Productivity (PROD)",
    "PROD estimate is normal (7)",
    "Justification: students",
    "Some code reuse is anticipated so adjusted effort = NOP",
    "NOP = (34 obj pts x (100-.1)) / 100 = 33.97",
    "Effort Estimate in person months (PM) (PM = 32)",
    "Effort(PM) = 33.97/7 = 4.85 person months",
    "Effort in Calendar Schedule",
    "Fisal @ 9 hours per week",
    "Kylie @ 9 hours per week",
    "Cory @ 9 hours per week",
    "Morgan @ 9 hours per week",
    "Cyler @ 9 hours per week",
    "Bijan @ 9 hours per week",
    "54 hours total per week",
    "216 hours (54x4) per month on the project",
    "3.59 months " + [char]0x2013 + " (160 hours/person month x 4.85 person months) / 216 hours/month"
)

$joined = [string]::Join("`r", $lines)
$tr9.Text = $joined

$levels = @(0,1,1,0,1,0,1,0,1,1,1,1,1,1,1,1,1)
for ($i = 0; $i -lt $levels.Length; $i++) {
    $para = $tr9.Paragraphs($i + 1, 1)
    $para.IndentLevel = $levels[$i] + 1
}
